# Update "想去人数" (F column) figures to the latest scraped counts.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1077
$ws1.Range("F4").Value = 1603
$ws1.Range("F5").Value = 730
$ws1.Range("F6").Value = 43

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 11

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1077
$ws4.Range("F4").Value = 1603
$ws4.Range("F5").Value = 11
$ws4.Range("F6").Value = 730
$ws4.Range("F7").Value = 43
